$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, pushing existing rows 4.. down by one.
$ws.Rows.Item(4).Insert()

# The former row 4 (now row 5) holds the template for the new row: duplicate it
# into the newly-inserted blank row 4, then adjust the two cells that differ.
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(4).PasteSpecial()
$excel.CutCopyMode = $false

# New row's own values (date + volume); price columns keep the copied value.
$ws.Cells.Item(4, 4).Value = 45043   # D4 - Fecha
$ws.Cells.Item(4, 10).Value = 120    # J4 - Volumen
